{"js": "// Add a new centered/bold/68-sz title line \"Untuk pemula\" right after the\n// existing \"TUTORIAL MEMBUAT GIT\" heading paragraph. The new line is typed\n// as two runs (\"Untuk\" / \" pemula\") wrapped with spell-check proofErr\n// markers around the first word, matching what Word itself emits, so we\n// build the paragraph from raw OOXML rather than insertParagraph/insertText.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst heading = paragraphs.items[0];\nconst afterHeading = heading.getRange(\"After\");\n\nconst newParagraphOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:pPr>' +\n                '<w:jc w:val=\"center\"/>' +\n                '<w:rPr><w:b/><w:sz w:val=\"68\"/></w:rPr>' +\n              '</w:pPr>' +\n              '<w:proofErr w:type=\"spellStart\"/>' +\n              '<w:r>' +\n                '<w:rPr><w:b/><w:sz w:val=\"68\"/></w:rPr>' +\n                '<w:t>Untuk</w:t>' +\n              '</w:r>' +\n              '<w:proofErr w:type=\"spellEnd\"/>' +\n              '<w:r>' +\n                '<w:rPr><w:b/><w:sz w:val=\"68\"/></w:rPr>' +\n                '<w:t xml:space=\"preserve\"> pemula</w:t>' +\n              '</w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\nafterHeading.insertOoxml(newParagraphOoxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Add a new centered/bold/68-sz title line \"Untuk pemula\" right after the\n# existing \"TUTORIAL MEMBUAT GIT\" heading. Word itself would type this as\n# two runs (\"Untuk\" / \" pemula\") wrapped with spell-check proofErr markers\n# around the first word, so build the paragraph from raw WordprocessingML\n# and insert it via Range.InsertXML instead of just setting .Text.\n\n$d = $word.ActiveDocument\n\n$r = $d.Content\n$r.Collapse(0)  # wdCollapseEnd -- land right after the last paragraph\n\n$xml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n          '<w:p>' +\n            '<w:pPr>' +\n              '<w:jc w:val=\"center\"/>' +\n              '<w:rPr><w:b/><w:sz w:val=\"68\"/></w:rPr>' +\n            '</w:pPr>' +\n            '<w:proofErr w:type=\"spellStart\"/>' +\n            '<w:r>' +\n              '<w:rPr><w:b/><w:sz w:val=\"68\"/></w:rPr>' +\n              '<w:t>Untuk</w:t>' +\n            '</w:r>' +\n            '<w:proofErr w:type=\"spellEnd\"/>' +\n            '<w:r>' +\n              '<w:rPr><w:b/><w:sz w:val=\"68\"/></w:rPr>' +\n              '<w:t xml:space=\"preserve\"> pemula</w:t>' +\n            '</w:r>' +\n          '</w:p>' +\n        '</w:body>' +\n      '</w:document>' +\n    '</pkg:xmlData>' +\n  '</pkg:part>' +\n'</pkg:package>'\n\n$r.InsertXML($xml)\n"}
